# order_20.xlsx -- "ready to plot distributions"
#
# 1. Row 20 held a bad data point (Stop time ~7h after Start instead of a
#    couple of minutes) -- drop it, which shifts the remaining rows up.
# 2. The trailing blank row (old row 28, now shifted to row 27) goes too.
# 3. Every Start/Stop serial value is nudged forward by exactly one day
#    (the date portion was off by a day; h:mm:ss display is unaffected).
# 4. Column B picks up the same left-aligned style column A already uses.
# 5. The "Stop" header font's color is pinned to explicit black instead of
#    the implicit theme color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- drop the outlier row and the trailing empty row -----------------
$ws.Rows("20").Delete()
$ws.Rows("27").Delete()

# --- shift every remaining timestamp forward by one day --------------
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 1).Value2 + 1.0
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 2).Value2 + 1.0
}

# --- column B now matches column A's (left-aligned) style ------------
$ws.Range("A2:A26").Copy()
$ws.Range("B2:B26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- "Stop" header: explicit black instead of theme color ------------
$ws.Range("B1").Font.Color = 0
